# Add new incident rows (35-45) to the management log, per the "new
# methods and fixing some incidences / v.2.0" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("WC47 NACP", "Fallo en elevador", "2024-05-30", "16:40:03", "Mañana", "16:40:07", "0:00:04", "N/A"),
    @("WC47 NACP", "Fallo en paletizador", "2024-05-30", "16:41:16", "Mañana", "16:41:20", "0:00:04", "N/A"),
    @("WC47 NACP", "Fallo tolva", "2024-05-30", "16:47:05", "Mañana", "16:47:27", "0:00:22", "N/A"),
    @("WC47 NACP", "Ascensor no sube", "2024-05-30", "16:47:09", "Mañana", "16:47:28", "0:00:19", "N/A"),
    @("WC47 NACP", "Fallo fijador tapa", "2024-05-30", "16:47:19", "Mañana", "16:47:32", "0:00:13", "0.07 minutos"),
    @("WC47 NACP", "No pone tornillo", "2024-05-30", "16:47:24", "Mañana", "16:47:30", "0:00:06", "0.11 minutos"),
    @("WC48 P5F", "Cámara no detecta Pcb", "2024-05-30", "17:02:51", "Mañana", "17:02:54", "0:00:03", "N/A"),
    @("WC48 P5F", "Detección de sealling mal puesto", "2024-05-30", "17:02:56", "Mañana", "17:03:00", "0:00:04", "N/A"),
    @("WC48 P5F", "Cámara no detecta skeleton", "2024-05-30", "17:57:13", "Mañana", "", "", "N/A"),
    @("WC47 NACP", "Fallo en elevador", "2024-05-30", "18:00:16", "Noche", "", "", "N/A"),
    @("WC48 P5F", "Cámara no detecta busbar", "2024-05-30", "18:21:34", "Noche", "", "", "N/A")
)

$startRow = 35
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]

    # The "Fecha" column holds plain text dates (e.g. "2024-05-30"), not
    # real date serials, matching how the rest of the sheet stores them.
    # A leading apostrophe forces Excel to keep the literal as text
    # instead of auto-converting it to a date value; resetting the style
    # afterwards drops the quote-prefix formatting flag it leaves behind.
    $ws.Cells.Item($row, 3).Value = "'" + $vals[2]
    $ws.Cells.Item($row, 3).Style = "Normal"

    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]

    # Rows 43/44/45 have no repair time recorded yet (still-open
    # incidents), so columns F/G are left blank for them.
    if ($vals[5] -ne "") {
        $ws.Cells.Item($row, 6).Value = $vals[5]
    }
    if ($vals[6] -ne "") {
        $ws.Cells.Item($row, 7).Value = $vals[6]
    }
    $ws.Cells.Item($row, 8).Value = $vals[7]
}
